# Apply data update to the "Inscricoes" worksheet
# Updates Inscritos (E), Pagos (F) and Inscricoes homologadas (H) columns
# for the rows affected by this data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$updates = @{
    "E3" = 11
    "F3" = 7
    "H3" = 7
    "E6" = 2
    "E8" = 6
    "F15" = 42
    "H15" = 42
    "E17" = 59
    "F17" = 24
    "H17" = 24
    "E18" = 50
    "F18" = 15
    "H18" = 15
    "F19" = 17
    "H19" = 17
    "E24" = 13
    "F24" = 6
    "H24" = 6
    "E26" = 13
    "F26" = 7
    "H26" = 7
    "E27" = 5
    "F27" = 2
    "H27" = 2
    "E33" = 17
    "E34" = 8
    "F34" = 2
    "H34" = 2
    "E36" = 49
    "F36" = 17
    "H36" = 17
    "E37" = 23
    "F38" = 8
    "H38" = 8
    "E39" = 12
    "F39" = 6
    "H39" = 6
    "E41" = 18
    "E42" = 16
    "F42" = 4
    "H42" = 4
    "F43" = 8
    "H43" = 8
    "E47" = 38
    "F47" = 20
    "H47" = 20
    "E49" = 34
    "F49" = 17
    "H49" = 17
    "E50" = 10
    "F50" = 1
    "H50" = 1
    "E60" = 10
    "F60" = 3
    "H60" = 3
    "E61" = 14
    "E63" = 12
    "E64" = 22
    "F65" = 5
    "H65" = 5
    "E66" = 22
    "E67" = 20
    "F67" = 9
    "H67" = 9
    "E70" = 17
    "E71" = 14
    "F71" = 10
    "H71" = 10
    "F72" = 9
    "H72" = 9
    "E73" = 15
    "E75" = 8
    "E76" = 27
    "F76" = 8
    "H76" = 8
    "E77" = 26
    "F77" = 7
    "H77" = 7
    "E79" = 13
    "E85" = 2
    "F85" = 1
    "H85" = 1
    "E89" = 13
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
